$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 21.18599966666667
$ws.Range("H2").Value = 63.557999
$ws.Range("I2").Value = 0.08765141600314529
$ws.Range("J2").Value = 0.08765141600314529
$ws.Range("M2").Value = 55.848606
$ws.Range("N2").Value = 167.545818
$ws.Range("O2").Value = 0.2323375192077237
$ws.Range("P2").Value = 0.2323375192077236
$ws.Range("Q2").Value = 1183.208548099798
$ws.Range("R2").Value = 10648.87693289818
$ws.Range("S2").Value = 0.02036471254921494
$ws.Range("T2").Value = 0.02036471254921494
$ws.Range("G3").Value = 21.18599966666667
$ws.Range("H3").Value = 63.557999
$ws.Range("I3").Value = 0.08765141600314529
$ws.Range("J3").Value = 0.08765141600314529
$ws.Range("O3").Value = 0.3515710112922583
$ws.Range("P3").Value = 0.3515710112922583
$ws.Range("Q3").Value = 1790.420364492133
$ws.Range("R3").Value = 16113.7832804292
$ws.Range("S3").Value = 0.03081569696542422
$ws.Range("T3").Value = 0.03081569696542422
$ws.Range("G4").Value = 21.18599966666667
$ws.Range("H4").Value = 63.557999
$ws.Range("I4").Value = 0.08765141600314529
$ws.Range("J4").Value = 0.08765141600314529
$ws.Range("M4").Value = 33.195992
$ws.Range("N4").Value = 99.58797600000001
$ws.Range("O4").Value = 0.1380996766314891
$ws.Range("P4").Value = 0.1380996766314891
$ws.Range("Q4").Value = 703.2902754466694
$ws.Range("R4").Value = 6329.612479020025
$ws.Range("S4").Value = 0.01210463220632649
$ws.Range("T4").Value = 0.01210463220632649
$ws.Range("G5").Value = 21.18599966666667
$ws.Range("H5").Value = 63.557999
$ws.Range("I5").Value = 0.08765141600314529
$ws.Range("J5").Value = 0.08765141600314529
$ws.Range("M5").Value = 66.82284533333335
$ws.Range("N5").Value = 200.468536
$ws.Range("O5").Value = 0.277991792868529
$ws.Range("P5").Value = 0.2779917928685289
$ws.Range("Q5").Value = 1415.708778957719
$ws.Range("R5").Value = 12741.37901061947
$ws.Range("S5").Value = 0.02436637428217963
$ws.Range("T5").Value = 0.02436637428217963
$ws.Range("I6").Value = 0.5040014103551328
$ws.Range("J6").Value = 0.5040014103551328
$ws.Range("M6").Value = 55.848606
$ws.Range("N6").Value = 167.545818
$ws.Range("O6").Value = 0.2323375192077237
$ws.Range("P6").Value = 0.2323375192077236
$ws.Range("Q6").Value = 6803.527018492754
$ws.Range("R6").Value = 61231.74316643479
$ws.Range("S6").Value = 0.1170984373591055
$ws.Range("T6").Value = 0.1170984373591055
$ws.Range("I7").Value = 0.5040014103551328
$ws.Range("J7").Value = 0.5040014103551328
$ws.Range("O7").Value = 0.3515710112922583
$ws.Range("P7").Value = 0.3515710112922583
$ws.Range("S7").Value = 0.1771922855312785
$ws.Range("T7").Value = 0.1771922855312785
$ws.Range("I8").Value = 0.5040014103551328
$ws.Range("J8").Value = 0.5040014103551328
$ws.Range("M8").Value = 33.195992
$ws.Range("N8").Value = 99.58797600000001
$ws.Range("O8").Value = 0.1380996766314891
$ws.Range("P8").Value = 0.1380996766314891
$ws.Range("Q8").Value = 4043.965367330195
$ws.Range("R8").Value = 36395.68830597176
$ws.Range("S8").Value = 0.06960243179185829
$ws.Range("T8").Value = 0.06960243179185828
$ws.Range("I9").Value = 0.5040014103551328
$ws.Range("J9").Value = 0.5040014103551328
$ws.Range("M9").Value = 66.82284533333335
$ws.Range("N9").Value = 200.468536
$ws.Range("O9").Value = 0.277991792868529
$ws.Range("P9").Value = 0.2779917928685289
$ws.Range("Q9").Value = 8140.418646758989
$ws.Range("R9").Value = 73263.76782083089
$ws.Range("S9").Value = 0.1401082556728906
$ws.Range("T9").Value = 0.1401082556728905
$ws.Range("G10").Value = 37.20718233333333
$ws.Range("H10").Value = 111.621547
$ws.Range("I10").Value = 0.1539347809079331
$ws.Range("J10").Value = 0.1539347809079331
$ws.Range("M10").Value = 55.848606
$ws.Range("N10").Value = 167.545818
$ws.Range("O10").Value = 0.2323375192077237
$ws.Range("P10").Value = 0.2323375192077236
$ws.Range("Q10").Value = 2077.969266504494
$ws.Range("R10").Value = 18701.72339854045
$ws.Range("S10").Value = 0.03576482511593364
$ws.Range("T10").Value = 0.03576482511593364
$ws.Range("G11").Value = 37.20718233333333
$ws.Range("H11").Value = 111.621547
$ws.Range("I11").Value = 0.1539347809079331
$ws.Range("J11").Value = 0.1539347809079331
$ws.Range("O11").Value = 0.3515710112922583
$ws.Range("P11").Value = 0.3515710112922583
$ws.Range("Q11").Value = 3144.36410851946
$ws.Range("R11").Value = 28299.27697667514
$ws.Range("S11").Value = 0.05411900659685425
$ws.Range("T11").Value = 0.05411900659685426
$ws.Range("G12").Value = 37.20718233333333
$ws.Range("H12").Value = 111.621547
$ws.Range("I12").Value = 0.1539347809079331
$ws.Range("J12").Value = 0.1539347809079331
$ws.Range("M12").Value = 33.195992
$ws.Range("N12").Value = 99.58797600000001
$ws.Range("O12").Value = 0.1380996766314891
$ws.Range("P12").Value = 0.1380996766314891
$ws.Range("Q12").Value = 1235.129327079875
$ws.Range("R12").Value = 11116.16394371887
$ws.Range("S12").Value = 0.02125834346572469
$ws.Range("T12").Value = 0.02125834346572469
$ws.Range("G13").Value = 37.20718233333333
$ws.Range("H13").Value = 111.621547
$ws.Range("I13").Value = 0.1539347809079331
$ws.Range("J13").Value = 0.1539347809079331
$ws.Range("M13").Value = 66.82284533333335
$ws.Range("N13").Value = 200.468536
$ws.Range("O13").Value = 0.277991792868529
$ws.Range("P13").Value = 0.2779917928685289
$ws.Range("Q13").Value = 2486.289790349466
$ws.Range("R13").Value = 22376.60811314519
$ws.Range("S13").Value = 0.04279260572942053
$ws.Range("T13").Value = 0.04279260572942053
$ws.Range("G14").Value = 61.49336899999999
$ws.Range("H14").Value = 184.480107
$ws.Range("I14").Value = 0.2544123927337887
$ws.Range("J14").Value = 0.2544123927337887
$ws.Range("M14").Value = 55.848606
$ws.Range("N14").Value = 167.545818
$ws.Range("O14").Value = 0.2323375192077237
$ws.Range("P14").Value = 0.2323375192077236
$ws.Range("Q14").Value = 3434.318936893614
$ws.Range("R14").Value = 30908.87043204252
$ws.Range("S14").Value = 0.05910954418346958
$ws.Range("T14").Value = 0.05910954418346957
$ws.Range("G15").Value = 61.49336899999999
$ws.Range("H15").Value = 184.480107
$ws.Range("I15").Value = 0.2544123927337887
$ws.Range("J15").Value = 0.2544123927337887
$ws.Range("O15").Value = 0.3515710112922583
$ws.Range("P15").Value = 0.3515710112922583
$ws.Range("Q15").Value = 5196.780037340189
$ws.Range("R15").Value = 46771.0203360617
$ws.Range("S15").Value = 0.08944402219870129
$ws.Range("T15").Value = 0.08944402219870129
$ws.Range("G16").Value = 61.49336899999999
$ws.Range("H16").Value = 184.480107
$ws.Range("I16").Value = 0.2544123927337887
$ws.Range("J16").Value = 0.2544123927337887
$ws.Range("M16").Value = 33.195992
$ws.Range("N16").Value = 99.58797600000001
$ws.Range("O16").Value = 0.1380996766314891
$ws.Range("P16").Value = 0.1380996766314891
$ws.Range("Q16").Value = 2041.333385377048
$ws.Range("R16").Value = 18372.00046839343
$ws.Range("S16").Value = 0.03513426916757963
$ws.Range("T16").Value = 0.03513426916757963
$ws.Range("G17").Value = 61.49336899999999
$ws.Range("H17").Value = 184.480107
$ws.Range("I17").Value = 0.2544123927337887
$ws.Range("J17").Value = 0.2544123927337887
$ws.Range("M17").Value = 66.82284533333335
$ws.Range("N17").Value = 200.468536
$ws.Range("O17").Value = 0.277991792868529
$ws.Range("P17").Value = 0.2779917928685289
$ws.Range("Q17").Value = 4109.161885712595
$ws.Range("R17").Value = 36982.45697141335
$ws.Range("S17").Value = 0.07072455718403825
$ws.Range("T17").Value = 0.07072455718403824
